# Weekly update: a new price record (row) is published for this
# market/product sub-set. It is inserted as the new first data row of the
# block (row 34), pushing all the existing rows in the block down by one.
#
# Net effect on the sheet: one row inserted at row 34 with brand-new data;
# every row that used to be 34..64 becomes 35..65 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 34 (shifts 34..64 -> 35..65).
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with this week's record.
$ws.Cells.Item(34, 1).Value  = 8
$ws.Cells.Item(34, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(34, 3).Value  = "Coquimbo"
$ws.Cells.Item(34, 4).Value  = 44893
$ws.Cells.Item(34, 5).Value  = 4
$ws.Cells.Item(34, 6).Value  = "Fruta"
$ws.Cells.Item(34, 7).Value  = 100103
$ws.Cells.Item(34, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(34, 9).Value  = 100103003
$ws.Cells.Item(34, 10).Value = "Damasco"
$ws.Cells.Item(34, 11).Value = "Castle Brite"
$ws.Cells.Item(34, 12).Value = "Primera"
$ws.Cells.Item(34, 13).Value = 160
$ws.Cells.Item(34, 14).Value = 26000
$ws.Cells.Item(34, 15).Value = 27000
$ws.Cells.Item(34, 16).Value = 26500
$ws.Cells.Item(34, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(34, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(34, 19).Value = 1656
$ws.Cells.Item(34, 20).Value = 16
